$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Value)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Value
    $Cell.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '37.807.71'
Set-TextValue $ws.Range('E2') '  +1.44%  '

Set-TextValue $ws.Range('D3') '2.107.05'
Set-TextValue $ws.Range('E3') '  +2.76%  '

Set-TextValue $ws.Range('E4') '  -0.05%  '

Set-TextValue $ws.Range('D5') '235.32'
Set-TextValue $ws.Range('E5') '  +1.18%  '

Set-TextValue $ws.Range('E6') '  +1.23%  '

Set-TextValue $ws.Range('D7') '58.38'
Set-TextValue $ws.Range('E7') '  +1.35%  '

Set-TextValue $ws.Range('E8') '  +0.00%  '

Set-TextValue $ws.Range('D9') '0.392'
Set-TextValue $ws.Range('E9') '  +2.36%  '

Set-TextValue $ws.Range('D10') '0.0780'
Set-TextValue $ws.Range('E10') '  +2.87%  '

Set-TextValue $ws.Range('E11') '  +1.28%  '

Set-TextValue $ws.Range('D12') '2.418.09'
Set-TextValue $ws.Range('E12') '  +2.64%  '

Set-TextValue $ws.Range('D13') '14.55'
Set-TextValue $ws.Range('E13') '  +1.58%  '

Set-TextValue $ws.Range('D14') '21.21'
Set-TextValue $ws.Range('E14') '  +1.86%  '

Set-TextValue $ws.Range('E15') '  +1.86%  '

Set-TextValue $ws.Range('D16') '5.24'
Set-TextValue $ws.Range('E16') '  +1.56%  '

Set-TextValue $ws.Range('D17') '2.110.70'
Set-TextValue $ws.Range('E17') '  +2.80%  '

Set-TextValue $ws.Range('D18') '37.759.77'
Set-TextValue $ws.Range('E18') '  +0.69%  '

Set-TextValue $ws.Range('D19') '6.24'
Set-TextValue $ws.Range('E19') '  +1.55%  '

Set-TextValue $ws.Range('D20') '70.12'
Set-TextValue $ws.Range('E20') '  +1.64%  '

Set-TextValue $ws.Range('D21') '0.0₃0822'
Set-TextValue $ws.Range('E21') '  +1.67%  '

Set-TextValue $ws.Range('D22') '227.85'
Set-TextValue $ws.Range('E22') '  +1.07%  '

Set-TextValue $ws.Range('E23') '  -0.02%  '

Set-TextValue $ws.Range('E24') '  +0.80%  '

Set-TextValue $ws.Range('E25') '  +0.39%  '

Set-TextValue $ws.Range('E26') '  +1.71%  '

Set-TextValue $ws.Range('E27') '  +1.28%  '

Set-TextValue $ws.Range('D28') '0.134'
Set-TextValue $ws.Range('E28') '  +3.41%  '

Set-TextValue $ws.Range('D29') '1.42'
Set-TextValue $ws.Range('E29') '  -3.90%  '

Set-TextValue $ws.Range('D30') '19.49'
Set-TextValue $ws.Range('E30') '  +2.21%  '

Set-TextValue $ws.Range('E31') '  +0.77%  '

Set-TextValue $ws.Range('D32') '4.66'
Set-TextValue $ws.Range('E32') '  +3.88%  '

Set-TextValue $ws.Range('D33') '2.59'
Set-TextValue $ws.Range('E33') '  +1.54%  '

Set-TextValue $ws.Range('D34') '0.0622'
Set-TextValue $ws.Range('E34') '  +0.51%  '

Set-TextValue $ws.Range('E35') '  +0.48%  '

Set-TextValue $ws.Range('E36') '  +6.01%  '

Set-TextValue $ws.Range('D37') '1.78'
Set-TextValue $ws.Range('E37') '  +1.01%  '

Set-TextValue $ws.Range('E38') '  -0.10%  '

Set-TextValue $ws.Range('D39') '5.62'
Set-TextValue $ws.Range('E39') '  -4.93%  '

Set-TextValue $ws.Range('E40') '  +0.00%  '

Set-TextValue $ws.Range('E41') '  +2.46%  '

Set-TextValue $ws.Range('D42') '97.60'
Set-TextValue $ws.Range('E42') '  +2.30%  '

Set-TextValue $ws.Range('D43') '1.475.53'
Set-TextValue $ws.Range('E43') '  +1.41%  '

Set-TextValue $ws.Range('E44') '  +1.48%  '

Set-TextValue $ws.Range('E45') '  -0.52%  '

Set-TextValue $ws.Range('D46') '4.21'
Set-TextValue $ws.Range('E46') '  -9.35%  '

Set-TextValue $ws.Range('B47') 'FraxShare'
Set-TextValue $ws.Range('C47') 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range('D47') '7.48'
Set-TextValue $ws.Range('E47') '  +4.55%  '

Set-TextValue $ws.Range('B48') 'ARBITRUM'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range('D48') '1.05'
Set-TextValue $ws.Range('E48') '  +2.62%  '

Set-TextValue $ws.Range('B49') 'InjectiveProtocol'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range('D49') '15.59'
Set-TextValue $ws.Range('E49') '  -0.15%  '

Set-TextValue $ws.Range('E50') '  +3.47%  '

Set-TextValue $ws.Range('D51') '2.304.63'
Set-TextValue $ws.Range('E51') '  +2.79%  '
